# Apply the "numbers suck but trying to fix them" edit:
#  - Remove the "flare" row (A8:C8) from the "IEAGHG chains" sheet and clear
#    the now-orphaned "heat_flare" value in E8 (keeping its style).
#  - Drop the text-format style that had been applied to E7.
#  - Update the selections/active sheet left behind by the edit session.

$wb = $excel.ActiveWorkbook

$wsChains      = $wb.Worksheets.Item("IEAGHG chains")
$wsConnections = $wb.Worksheets.Item("IEAGHG connections")
$wsSteelChain  = $wb.Worksheets.Item("IEAGHG steel chain")

# --- "IEAGHG chains": delete the flare/heat/inflow row, tidy up E7/E8 ---
$wsChains.Range("A8:C8").ClearContents()
$wsChains.Range("E8").ClearContents()
$wsChains.Range("E7").Style = "Normal"

# --- selections left over from the editing session ---
$wsConnections.Range("B17:C17").Select()

$wsSteelChain.Range("C4").Select()

# "IEAGHG chains" ends up the active sheet/selection when the file was saved
$wsChains.Activate()
$wsChains.Range("B14").Select()
